$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Insert a new row above the current row 8 ("AccountDetails"/Magento block),
# shifting all rows 8-35 down to 9-36.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the "Invoice" entry.
$ws.Cells.Item(8, 1).Value = "Invoice"
$ws.Cells.Item(8, 7).Value = "Invoice"
$ws.Cells.Item(8, 8).Value = "https://na-preprod.hele.digital/rest/ospreyusen/V1/order/"

# Widen column H so the new long URL fits (close to the 72.140625 bestFit width).
$ws.Columns.Item(8).ColumnWidth = 71.3

# This runtime does not shift hyperlink ranges when rows are inserted, so we
# rebuild the whole hyperlinks collection to match the post-insert layout
# (existing links shifted down one row where needed, plus the new H8 link).
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:qatesting.lotuswave@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:Paypal-buyer@hydroflask.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:gsapram@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C32"), "mailto:qatesting.lotuswave@gmail.com", "", "", "mailto:qatesting.lotuswave@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E36"), "mailto:Lotuswave@123", "", "", "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D36"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C36"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E35"), "mailto:Lotuswave@123", "", "", "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D35"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C35"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B35"), "mailto:skatipelli@helenoftroy.com", "", "", "mailto:skatipelli@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:testersemail.278@gmail.com", "", "", "testersemail.278@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:lotusqa.gld.stg.os.us.automation.01@gmail.com", "", "", "mailto:lotusqa.gld.stg.os.us.automation.01@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B21"), "mailto:lotusqa.gld.stg.os.us.automation.01@gmail.com", "", "", "mailto:lotusqa.gld.stg.os.us.automation.01@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E22"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "mailto:Lotuswave@1234") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), "mailto:testersemail.278@gmail.com", "", "", "testersemail.278@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:Paypal-buyer@hydroflask.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:testersemail.278@gmail.com", "", "", "testersemail.278@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:rrendla@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J5"), "https://na-preprod.hele.digital/rest/V1/integration/admin/token", "", "", "https://na-preprod.hele.digital/rest/V1/integration/admin/token") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:avayugundla@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H9"), "https://na-preprod.hele.digital/heledigitaladmin/admin/dashboard/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H10"), "https://na-preprod.hele.digital/heledigitaladmin/sales/order/view/order_id/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H6"), "https://na-preprod.hele.digital/rest/ospreyusen/V1/orders/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H7"), "https://na-preprod.hele.digital/rest/all/V1/order/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H8"), "https://na-preprod.hele.digital/rest/ospreyusen/V1/order/") | Out-Null

# Restore the active cell selection to match the edited workbook.
$ws.Range("H19").Select()
